$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifting dbExcel/WebExcel columns right)
$ws.Columns("B").Insert()

# Set header row
$ws.Range("B1").Value = "StatQuery"

# Set new query cell with wrap text style like A2
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.ethnicity IN ['UNKNOWN'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Ensure column B picks up the same width as column A
$ws.Columns("B").ColumnWidth = 75.81640625
